# Auto-applies the numeric cell updates captured in the commit diff.
# Values are written via .Value2 (plain numeric) on Cells.Item(row, col);
# cells that must become empty use ClearContents(); cells that gain a
# value for the first time are simply assigned.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value2 = 691.53845  # H19
$ws.Cells.Item(19, 9).Value2 = 832  # I19
$ws.Cells.Item(19, 11).Value2 = 832  # K19
$ws.Cells.Item(19, 13).Value2 = -657  # M19
$ws.Cells.Item(64, 8).Value2 = 7986.9414  # H64
$ws.Cells.Item(64, 9).Value2 = 6597.25  # I64
$ws.Cells.Item(64, 11).Value2 = 6597.25  # K64
$ws.Cells.Item(64, 13).Value2 = -6349.25  # M64
$ws.Cells.Item(67, 8).Value2 = 7986.9414  # H67
$ws.Cells.Item(67, 9).Value2 = 6597.25  # I67
$ws.Cells.Item(67, 11).Value2 = 6597.25  # K67
$ws.Cells.Item(67, 13).Value2 = -5739.25  # M67
$ws.Cells.Item(106, 8).Value2 = 2520.7  # H106
$ws.Cells.Item(106, 9).Value2 = 2520.7  # I106
$ws.Cells.Item(106, 11).Value2 = 2520.7  # K106
$ws.Cells.Item(106, 13).Value2 = -1889.7  # M106
$ws.Cells.Item(127, 8).Value2 = 9338.974  # H127
$ws.Cells.Item(127, 9).Value2 = 2094.8635  # I127
$ws.Cells.Item(127, 10).Value2 = 19299.625  # J127
$ws.Cells.Item(127, 11).Value2 = 6284.5905  # K127
$ws.Cells.Item(127, 12).Value2 = 57898.875  # L127
$ws.Cells.Item(127, 13).Value2 = -1324.5905  # M127
$ws.Cells.Item(127, 14).Value2 = -67818.875  # N127
$ws.Cells.Item(132, 8).Value2 = 688.9138  # H132
$ws.Cells.Item(132, 9).Value2 = 697.5  # I132
$ws.Cells.Item(132, 11).Value2 = 2092.5  # K132
$ws.Cells.Item(132, 13).Value2 = 437.5  # M132
$ws.Cells.Item(135, 8).Value2 = 5929.1665  # H135
$ws.Cells.Item(135, 9).Value2 = 1422.6666  # I135
$ws.Cells.Item(135, 11).Value2 = 12803.9994  # K135
$ws.Cells.Item(135, 13).Value2 = -10268.9994  # M135
$ws.Cells.Item(137, 8).Value2 = 1839.9818  # H137
$ws.Cells.Item(137, 9).Value2 = 1355.1177  # I137
$ws.Cells.Item(137, 11).Value2 = 4065.3531  # K137
$ws.Cells.Item(137, 13).Value2 = -1515.3531  # M137
$ws.Cells.Item(138, 8).Value2 = 3878.4688  # H138
$ws.Cells.Item(138, 9).Value2 = 1999.3334  # I138
$ws.Cells.Item(138, 11).Value2 = 5998.0002  # K138
$ws.Cells.Item(138, 13).Value2 = -858.0002000000004  # M138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 17808.186  # H32
$ws.Cells.Item(32, 9).Value2 = 13774.869  # I32
$ws.Cells.Item(32, 11).Value2 = 13774.869  # K32
$ws.Cells.Item(32, 13).Value2 = -13487.869  # M32
$ws.Cells.Item(45, 8).Value2 = 5597.7856  # H45
$ws.Cells.Item(45, 9).Value2 = 3047.182  # I45
$ws.Cells.Item(45, 11).Value2 = 3047.182  # K45
$ws.Cells.Item(45, 13).Value2 = -2670.182  # M45
$ws.Cells.Item(61, 8).Value2 = 4056.7144  # H61
$ws.Cells.Item(61, 9).Value2 = 1937.8462  # I61
$ws.Cells.Item(61, 11).Value2 = 1937.8462  # K61
$ws.Cells.Item(61, 13).Value2 = -1725.8462  # M61
$ws.Cells.Item(74, 8).Value2 = 5674.75  # H74
$ws.Cells.Item(74, 9).Value2 = 5950  # I74
$ws.Cells.Item(74, 11).Value2 = 5950  # K74
$ws.Cells.Item(74, 13).Value2 = -5076  # M74
$ws.Cells.Item(77, 8).Value2 = 5674.75  # H77
$ws.Cells.Item(77, 9).Value2 = 5950  # I77
$ws.Cells.Item(77, 11).Value2 = 29750  # K77
$ws.Cells.Item(77, 13).Value2 = -25382  # M77
$ws.Cells.Item(88, 8).Value2 = 6420194  # H88
$ws.Cells.Item(88, 9).Value2 = 19415.834  # I88
$ws.Cells.Item(88, 10).Value2 = 11906575  # J88
$ws.Cells.Item(88, 11).Value2 = 19415.834  # K88
$ws.Cells.Item(88, 12).Value2 = 11906575  # L88
$ws.Cells.Item(88, 13).Value2 = -19009.834  # M88
$ws.Cells.Item(88, 14).Value2 = -11907387  # N88
$ws.Cells.Item(91, 8).Value2 = 6420194  # H91
$ws.Cells.Item(91, 9).Value2 = 19415.834  # I91
$ws.Cells.Item(91, 10).Value2 = 11906575  # J91
$ws.Cells.Item(91, 11).Value2 = 19415.834  # K91
$ws.Cells.Item(91, 12).Value2 = 11906575  # L91
$ws.Cells.Item(91, 13).Value2 = -18011.834  # M91
$ws.Cells.Item(91, 14).Value2 = -11909383  # N91
$ws.Cells.Item(111, 8).Value2 = 0  # H111
$ws.Cells.Item(111, 10).Value2 = 0  # J111
$ws.Cells.Item(111, 12).Value2 = 0  # L111
$ws.Cells.Item(111, 14).ClearContents()  # N111
$ws.Cells.Item(132, 8).Value2 = 4527.641  # H132
$ws.Cells.Item(132, 9).Value2 = 4237  # I132
$ws.Cells.Item(132, 11).Value2 = 12711  # K132
$ws.Cells.Item(132, 13).Value2 = -10181  # M132
$ws.Cells.Item(136, 8).Value2 = 4056.7144  # H136
$ws.Cells.Item(136, 9).Value2 = 1937.8462  # I136
$ws.Cells.Item(136, 11).Value2 = 5813.5386  # K136
$ws.Cells.Item(136, 13).Value2 = -3263.5386  # M136
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value2 = 2379.8  # H105
$ws.Cells.Item(105, 9).Value2 = 1755.4166  # I105
$ws.Cells.Item(105, 11).Value2 = 1755.4166  # K105
$ws.Cells.Item(105, 13).Value2 = -8.416600000000017  # M105
$ws.Cells.Item(134, 8).Value2 = 2707.8572  # H134
$ws.Cells.Item(134, 9).Value2 = 1763.7428  # I134
$ws.Cells.Item(134, 11).Value2 = 5291.2284  # K134
$ws.Cells.Item(134, 13).Value2 = -2756.2284  # M134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 3363.6128  # H31
$ws.Cells.Item(31, 9).Value2 = 2380.4707  # I31
$ws.Cells.Item(31, 11).Value2 = 2380.4707  # K31
$ws.Cells.Item(31, 13).Value2 = -2085.4707  # M31
$ws.Cells.Item(34, 8).Value2 = 3363.6128  # H34
$ws.Cells.Item(34, 9).Value2 = 2380.4707  # I34
$ws.Cells.Item(34, 11).Value2 = 2380.4707  # K34
$ws.Cells.Item(34, 13).Value2 = -2178.4707  # M34
$ws.Cells.Item(99, 8).Value2 = 11224253  # H99
$ws.Cells.Item(99, 9).Value2 = 2714429  # I99
$ws.Cells.Item(99, 11).Value2 = 2714429  # K99
$ws.Cells.Item(99, 13).Value2 = -2712931  # M99
$ws.Cells.Item(126, 8).Value2 = 11224253  # H126
$ws.Cells.Item(126, 9).Value2 = 2714429  # I126
$ws.Cells.Item(126, 11).Value2 = 8143287  # K126
$ws.Cells.Item(126, 13).Value2 = -8140817  # M126
$ws.Cells.Item(132, 8).Value2 = 2248.9688  # H132
$ws.Cells.Item(132, 9).Value2 = 1865.5667  # I132
$ws.Cells.Item(132, 10).Value2 = 8000  # J132
$ws.Cells.Item(132, 11).Value2 = 5596.7001  # K132
$ws.Cells.Item(132, 12).Value2 = 24000  # L132
$ws.Cells.Item(132, 13).Value2 = -3066.7001  # M132
$ws.Cells.Item(132, 14).Value2 = -29060  # N132
$ws.Cells.Item(140, 8).Value2 = 71285.64  # H140
$ws.Cells.Item(140, 10).Value2 = 72153.84  # J140
$ws.Cells.Item(140, 12).Value2 = 72153.84  # L140
$ws.Cells.Item(140, 14).Value2 = -82513.84  # N140
$ws.Cells.Item(141, 8).Value2 = 387914.47  # H141
$ws.Cells.Item(141, 10).Value2 = 387914.47  # J141
$ws.Cells.Item(141, 12).Value2 = 387914.47  # L141
$ws.Cells.Item(141, 14).Value2 = -398274.47  # N141
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value2 = 4327.6665  # H5
$ws.Cells.Item(5, 9).Value2 = 337.5  # I5
$ws.Cells.Item(5, 11).Value2 = 1012.5  # K5
$ws.Cells.Item(5, 13).Value2 = -900.5  # M5
$ws.Cells.Item(99, 8).Value2 = 88399.664  # H99
$ws.Cells.Item(99, 9).Value2 = 750  # I99
$ws.Cells.Item(99, 10).Value2 = 132224.5  # J99
$ws.Cells.Item(99, 11).Value2 = 2250  # K99
$ws.Cells.Item(99, 12).Value2 = 396673.5  # L99
$ws.Cells.Item(99, 13).Value2 = -4  # M99
$ws.Cells.Item(99, 14).Value2 = -401165.5  # N99
$ws.Cells.Item(109, 8).Value2 = 2333.3333  # H109
$ws.Cells.Item(109, 9).Value2 = 500  # I109
$ws.Cells.Item(109, 10).Value2 = 3250  # J109
$ws.Cells.Item(109, 11).Value2 = 1500  # K109
$ws.Cells.Item(109, 12).Value2 = 9750  # L109
$ws.Cells.Item(109, 13).Value2 = -460  # M109
$ws.Cells.Item(109, 14).Value2 = -11830  # N109
$ws.Cells.Item(131, 8).Value2 = 4273.263  # H131
$ws.Cells.Item(131, 9).Value2 = 4209.875  # I131
$ws.Cells.Item(131, 10).Value2 = 4319.364  # J131
$ws.Cells.Item(131, 11).Value2 = 12629.625  # K131
$ws.Cells.Item(131, 12).Value2 = 12958.092  # L131
$ws.Cells.Item(131, 13).Value2 = -7589.625  # M131
$ws.Cells.Item(131, 14).Value2 = -23038.092  # N131
$ws.Cells.Item(135, 8).Value2 = 4327.6665  # H135
$ws.Cells.Item(135, 9).Value2 = 337.5  # I135
$ws.Cells.Item(135, 11).Value2 = 3037.5  # K135
$ws.Cells.Item(135, 13).Value2 = -502.5  # M135
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value2 = 0  # H68
$ws.Cells.Item(68, 10).Value2 = 0  # J68
$ws.Cells.Item(68, 12).Value2 = 0  # L68
$ws.Cells.Item(68, 14).ClearContents()  # N68
$ws.Cells.Item(71, 8).Value2 = 0  # H71
$ws.Cells.Item(71, 10).Value2 = 0  # J71
$ws.Cells.Item(71, 12).Value2 = 0  # L71
$ws.Cells.Item(71, 14).ClearContents()  # N71
$ws.Cells.Item(118, 8).Value2 = 35714.145  # H118
$ws.Cells.Item(118, 10).Value2 = 35714.145  # J118
$ws.Cells.Item(118, 12).Value2 = 35714.145  # L118
$ws.Cells.Item(118, 14).Value2 = -39028.145  # N118
$ws.Cells.Item(132, 8).Value2 = 4946.6924  # H132
$ws.Cells.Item(132, 9).Value2 = 4600.8096  # I132
$ws.Cells.Item(132, 10).Value2 = 6399.4  # J132
$ws.Cells.Item(132, 11).Value2 = 13802.4288  # K132
$ws.Cells.Item(132, 12).Value2 = 19198.2  # L132
$ws.Cells.Item(132, 13).Value2 = -11272.4288  # M132
$ws.Cells.Item(132, 14).Value2 = -24258.2  # N132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value2 = 5159.7334  # H68
$ws.Cells.Item(68, 9).Value2 = 3336.4736  # I68
$ws.Cells.Item(68, 11).Value2 = 3336.4736  # K68
$ws.Cells.Item(68, 13).Value2 = -2587.4736  # M68
$ws.Cells.Item(71, 8).Value2 = 5159.7334  # H71
$ws.Cells.Item(71, 9).Value2 = 3336.4736  # I71
$ws.Cells.Item(71, 11).Value2 = 16682.368  # K71
$ws.Cells.Item(71, 13).Value2 = -12938.368  # M71
$ws.Cells.Item(122, 8).Value2 = 5838.722  # H122
$ws.Cells.Item(122, 9).Value2 = 4401.6  # I122
$ws.Cells.Item(122, 11).Value2 = 13204.8  # K122
$ws.Cells.Item(122, 13).Value2 = -10754.8  # M122
$ws.Cells.Item(123, 8).Value2 = 60000  # H123
$ws.Cells.Item(123, 10).Value2 = 60000  # J123
$ws.Cells.Item(123, 12).Value2 = 60000  # L123
$ws.Cells.Item(123, 14).Value2 = -69800  # N123
$ws.Cells.Item(127, 8).Value2 = 0  # H127
$ws.Cells.Item(127, 10).Value2 = 0  # J127
$ws.Cells.Item(127, 12).Value2 = 0  # L127
$ws.Cells.Item(127, 14).ClearContents()  # N127
$ws.Cells.Item(132, 8).Value2 = 4297.1465  # H132
$ws.Cells.Item(132, 9).Value2 = 3738.7878  # I132
$ws.Cells.Item(132, 10).Value2 = 6600.375  # J132
$ws.Cells.Item(132, 11).Value2 = 11216.3634  # K132
$ws.Cells.Item(132, 12).Value2 = 19801.125  # L132
$ws.Cells.Item(132, 13).Value2 = -8686.3634  # M132
$ws.Cells.Item(132, 14).Value2 = -24861.125  # N132
$ws.Cells.Item(136, 8).Value2 = 4976.25  # H136
$ws.Cells.Item(136, 9).Value2 = 2931.3333  # I136
$ws.Cells.Item(136, 11).Value2 = 8793.999899999999  # K136
$ws.Cells.Item(136, 13).Value2 = -6243.999899999999  # M136
$ws.Cells.Item(137, 8).Value2 = 68333.25  # H137
$ws.Cells.Item(137, 9).Value2 = 59999  # I137
$ws.Cells.Item(137, 11).Value2 = 59999  # K137
$ws.Cells.Item(137, 13).Value2 = -54899  # M137
$ws.Cells.Item(139, 8).Value2 = 99166  # H139
$ws.Cells.Item(139, 9).Value2 = 99992  # I139
$ws.Cells.Item(139, 10).Value2 = 99090.91  # J139
$ws.Cells.Item(139, 11).Value2 = 99992  # K139
$ws.Cells.Item(139, 12).Value2 = 99090.91  # L139
$ws.Cells.Item(139, 13).Value2 = -94852  # M139
$ws.Cells.Item(139, 14).Value2 = -109370.91  # N139
$ws.Cells.Item(141, 8).Value2 = 80000  # H141
$ws.Cells.Item(141, 9).Value2 = 0  # I141
$ws.Cells.Item(141, 11).Value2 = 0  # K141
$ws.Cells.Item(141, 13).ClearContents()  # M141
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 11114440  # H62
$ws.Cells.Item(62, 10).Value2 = 13892363  # J62
$ws.Cells.Item(62, 12).Value2 = 13892363  # L62
$ws.Cells.Item(62, 14).Value2 = -13893611  # N62
$ws.Cells.Item(65, 8).Value2 = 11114440  # H65
$ws.Cells.Item(65, 10).Value2 = 13892363  # J65
$ws.Cells.Item(65, 12).Value2 = 69461815  # L65
$ws.Cells.Item(65, 14).Value2 = -69468055  # N65
$ws.Cells.Item(122, 8).Value2 = 4910.9443  # H122
$ws.Cells.Item(122, 9).Value2 = 3876.3076  # I122
$ws.Cells.Item(122, 11).Value2 = 11628.9228  # K122
$ws.Cells.Item(122, 13).Value2 = -9178.9228  # M122
$ws.Cells.Item(126, 8).Value2 = 1896.3  # H126
$ws.Cells.Item(126, 9).Value2 = 1551.4445  # I126
$ws.Cells.Item(126, 11).Value2 = 4654.333500000001  # K126
$ws.Cells.Item(126, 13).Value2 = -2184.333500000001  # M126
$ws.Cells.Item(128, 8).Value2 = 60536.25  # H128
$ws.Cells.Item(128, 10).Value2 = 60536.25  # J128
$ws.Cells.Item(128, 12).Value2 = 60536.25  # L128
$ws.Cells.Item(128, 14).Value2 = -70496.25  # N128
